# Fill in the results for the Mon, Feb 12, 2024 games (Sheet1 rows 178-187).
# Before this edit only the schedule (Date, Start, Away team, Home team, Arena)
# was known for these rows; this edit adds the actual scores, overtime flag,
# win/loss teams, forecasted pick, correctness, point-diff and validity.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Row 178: Indiana Pacers @ Charlotte Hornets ---
$ws1.Cells.Item(178,4).Value = 102
$ws1.Cells.Item(178,4).NumberFormat = "#,##0"
$ws1.Cells.Item(178,6).Value = 111
$ws1.Cells.Item(178,6).NumberFormat = "#,##0"
$ws1.Cells.Item(178,7).Value = "NA"
$ws1.Cells.Item(178,9).Value = "Charlotte Hornets"
$ws1.Cells.Item(178,10).Value = "Indiana Pacers"
$ws1.Cells.Item(178,11).Value = "Indiana Pacers"
$ws1.Cells.Item(178,12).Value = "No"
$ws1.Cells.Item(178,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(178,13).Formula = "=ABS(D178-F178)"
$ws1.Cells.Item(178,13).NumberFormat = "#,##0"
$ws1.Cells.Item(178,14).Formula = "=K178=I178"

# --- Row 179: Philadelphia 76ers @ Cleveland Cavaliers ---
$ws1.Cells.Item(179,4).Value = 123
$ws1.Cells.Item(179,4).NumberFormat = "#,##0"
$ws1.Cells.Item(179,6).Value = 121
$ws1.Cells.Item(179,6).NumberFormat = "#,##0"
$ws1.Cells.Item(179,7).Value = "NA"
$ws1.Cells.Item(179,9).Value = "Philadelphia 76ers"
$ws1.Cells.Item(179,10).Value = "Cleveland Cavaliers"
$ws1.Cells.Item(179,11).Value = "Cleveland Cavaliers"
$ws1.Cells.Item(179,12).Value = "No"
$ws1.Cells.Item(179,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(179,13).Formula = "=ABS(D179-F179)"
$ws1.Cells.Item(179,13).NumberFormat = "#,##0"
$ws1.Cells.Item(179,14).Formula = "=K179=I179"

# --- Row 180: Atlanta Hawks @ Chicago Bulls ---
$ws1.Cells.Item(180,4).Value = 136
$ws1.Cells.Item(180,4).NumberFormat = "#,##0"
$ws1.Cells.Item(180,6).Value = 126
$ws1.Cells.Item(180,6).NumberFormat = "#,##0"
$ws1.Cells.Item(180,7).Value = "NA"
$ws1.Cells.Item(180,9).Value = "Chicago Bulls"
$ws1.Cells.Item(180,10).Value = "Atlanta Hawks"
$ws1.Cells.Item(180,11).Value = "Chicago Bulls"
$ws1.Cells.Item(180,11).Interior.Color = 5287936
$ws1.Cells.Item(180,12).Value = "Yes"
$ws1.Cells.Item(180,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(180,13).Formula = "=ABS(D180-F180)"
$ws1.Cells.Item(180,13).NumberFormat = "#,##0"
$ws1.Cells.Item(180,14).Formula = "=K180=I180"

# --- Row 181: San Antonio Spurs @ Toronto Raptors ---
$ws1.Cells.Item(181,4).Value = 122
$ws1.Cells.Item(181,4).NumberFormat = "#,##0"
$ws1.Cells.Item(181,6).Value = 99
$ws1.Cells.Item(181,6).NumberFormat = "#,##0"
$ws1.Cells.Item(181,7).Value = "NA"
$ws1.Cells.Item(181,9).Value = "San Antonio Spurs"
$ws1.Cells.Item(181,10).Value = "Toronto Raptors"
$ws1.Cells.Item(181,11).Value = "Toronto Raptors"
$ws1.Cells.Item(181,12).Value = "No"
$ws1.Cells.Item(181,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(181,13).Formula = "=ABS(D181-F181)"
$ws1.Cells.Item(181,13).NumberFormat = "#,##0"
$ws1.Cells.Item(181,14).Formula = "=K181=I181"

# --- Row 182: New York Knicks @ Houston Rockets ---
$ws1.Cells.Item(182,4).Value = 103
$ws1.Cells.Item(182,4).NumberFormat = "#,##0"
$ws1.Cells.Item(182,6).Value = 105
$ws1.Cells.Item(182,6).NumberFormat = "#,##0"
$ws1.Cells.Item(182,7).Value = "NA"
$ws1.Cells.Item(182,9).Value = "Houston Rockets"
$ws1.Cells.Item(182,10).Value = "New York Knicks"
$ws1.Cells.Item(182,11).Value = "New York Knicks"
$ws1.Cells.Item(182,12).Value = "No"
$ws1.Cells.Item(182,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(182,13).Formula = "=ABS(D182-F182)"
$ws1.Cells.Item(182,13).NumberFormat = "#,##0"
$ws1.Cells.Item(182,14).Formula = "=K182=I182"

# --- Row 183: New Orleans Pelicans @ Memphis Grizzlies ---
$ws1.Cells.Item(183,4).Value = 96
$ws1.Cells.Item(183,4).NumberFormat = "#,##0"
$ws1.Cells.Item(183,6).Value = 87
$ws1.Cells.Item(183,6).NumberFormat = "#,##0"
$ws1.Cells.Item(183,7).Value = "NA"
$ws1.Cells.Item(183,9).Value = "New Orleans Pelicans"
$ws1.Cells.Item(183,10).Value = "Memphis Grizzlies"
$ws1.Cells.Item(183,11).Value = "New Orleans Pelicans"
$ws1.Cells.Item(183,11).Interior.Color = 5287936
$ws1.Cells.Item(183,12).Value = "Yes"
$ws1.Cells.Item(183,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(183,13).Formula = "=ABS(D183-F183)"
$ws1.Cells.Item(183,13).NumberFormat = "#,##0"
$ws1.Cells.Item(183,14).Formula = "=K183=I183"

# --- Row 184: Denver Nuggets @ Milwaukee Bucks ---
$ws1.Cells.Item(184,4).Value = 95
$ws1.Cells.Item(184,4).NumberFormat = "#,##0"
$ws1.Cells.Item(184,6).Value = 112
$ws1.Cells.Item(184,6).NumberFormat = "#,##0"
$ws1.Cells.Item(184,7).Value = "NA"
$ws1.Cells.Item(184,9).Value = "Milwaukee Bucks"
$ws1.Cells.Item(184,10).Value = "Denver Nuggets"
$ws1.Cells.Item(184,11).Value = "Denver Nuggets"
$ws1.Cells.Item(184,12).Value = "No"
$ws1.Cells.Item(184,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(184,13).Formula = "=ABS(D184-F184)"
$ws1.Cells.Item(184,13).NumberFormat = "#,##0"
$ws1.Cells.Item(184,14).Formula = "=K184=I184"

# --- Row 185: Washington Wizards @ Dallas Mavericks ---
$ws1.Cells.Item(185,4).Value = 104
$ws1.Cells.Item(185,4).NumberFormat = "#,##0"
$ws1.Cells.Item(185,6).Value = 112
$ws1.Cells.Item(185,6).NumberFormat = "#,##0"
$ws1.Cells.Item(185,7).Value = "NA"
$ws1.Cells.Item(185,9).Value = "Dallas Mavericks"
$ws1.Cells.Item(185,10).Value = "Washington Wizards"
$ws1.Cells.Item(185,11).Value = "Dallas Mavericks"
$ws1.Cells.Item(185,11).Interior.Color = 5287936
$ws1.Cells.Item(185,12).Value = "Yes"
$ws1.Cells.Item(185,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(185,13).Formula = "=ABS(D185-F185)"
$ws1.Cells.Item(185,13).NumberFormat = "#,##0"
$ws1.Cells.Item(185,14).Formula = "=K185=I185"

# --- Row 186: Golden State Warriors @ Utah Jazz ---
$ws1.Cells.Item(186,4).Value = 129
$ws1.Cells.Item(186,4).NumberFormat = "#,##0"
$ws1.Cells.Item(186,6).Value = 107
$ws1.Cells.Item(186,6).NumberFormat = "#,##0"
$ws1.Cells.Item(186,7).Value = "NA"
$ws1.Cells.Item(186,9).Value = "Golden State Warriors"
$ws1.Cells.Item(186,10).Value = "Utah Jazz"
$ws1.Cells.Item(186,11).Value = "Utah Jazz"
$ws1.Cells.Item(186,12).Value = "No"
$ws1.Cells.Item(186,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(186,13).Formula = "=ABS(D186-F186)"
$ws1.Cells.Item(186,13).NumberFormat = "#,##0"
$ws1.Cells.Item(186,14).Formula = "=K186=I186"

# --- Row 187: Los Angeles Clippers @ Minnesota Timberwolves ---
$ws1.Cells.Item(187,4).Value = 121
$ws1.Cells.Item(187,4).NumberFormat = "#,##0"
$ws1.Cells.Item(187,6).Value = 100
$ws1.Cells.Item(187,6).NumberFormat = "#,##0"
$ws1.Cells.Item(187,7).Value = "NA"
$ws1.Cells.Item(187,9).Value = "Minnesota Timberwolves"
$ws1.Cells.Item(187,10).Value = "Los Angeles Clippers"
$ws1.Cells.Item(187,11).Value = "Los Angeles Clippers"
$ws1.Cells.Item(187,12).Value = "No"
$ws1.Cells.Item(187,12).NumberFormat = "#,##0.00"
$ws1.Cells.Item(187,13).Formula = "=ABS(D187-F187)"
$ws1.Cells.Item(187,13).NumberFormat = "#,##0"
$ws1.Cells.Item(187,14).Formula = "=K187=I187"

# Recalculate so Sheet2's COUNTIFS / AVERAGEIFS summary rows (which read
# Sheet1!L:L, Sheet1!A:A, Sheet1!M:M) and the three trend charts that cache
# Sheet2!A2:A26 / D2:D26 / F2:F26 / G2:G26 pick up the new Feb 12 games.
$excel.Calculate()

# Restore the UI selection: the user ended up with Sheet2!P26 selected while
# Sheet1 remains the tab shown on reopen.
$ws2.Activate()
$ws2.Range("P26").Select()
$ws1.Activate()
